# Workbook edit: restructure the glossary workbook into two sheets
# ("Glossary" + "References") and add a reference/source link.

$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet to "Glossary" and update its selection ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Glossary"
$ws1.Range("B20").Select() | Out-Null

# --- Add the new "References" sheet right after "Glossary" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "References"

# --- Populate the References sheet ---
$ws2.Range("A1").Value = "llun (clystyrau_k_modd)"
$ws2.Range("B1").Value = "https://h1ros.github.io/posts/k-means-clustering/"

# Turn B1 into a real hyperlink (this also applies the built-in "Hyperlink" style)
$ws2.Hyperlinks.Add($ws2.Range("B1"), "https://h1ros.github.io/posts/k-means-clustering/") | Out-Null

# Match the column widths used for the References sheet
$ws2.Columns.Item(1).ColumnWidth = 26.3
$ws2.Columns.Item(2).ColumnWidth = 46.0

# Final selection on the References sheet (becomes the active/visible tab)
$ws2.Range("A2").Select() | Out-Null

Write-Host "Workbook restructured: Glossary + References sheets ready."
